$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Absent" column (H) to 1 for rows 3 through 18
$ws.Range("H3:H18").Value = 1
